# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record (Albahaca - Femacal de La Calera) was inserted
# as the new row 18 of the data table, pushing every existing row from 18
# downward down by one (old row 18 becomes 19, ..., old row 110 becomes 111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; Excel shifts rows 18..110 to 19..111
# and carries the date-format style from the surrounding rows onto the new row.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Cells.Item(18, 1).Value  = 3
$ws.Cells.Item(18, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(18, 3).Value  = "Coquimbo"
$ws.Cells.Item(18, 4).Value  = 44561
$ws.Cells.Item(18, 5).Value  = 5
$ws.Cells.Item(18, 6).Value  = 100112052
$ws.Cells.Item(18, 7).Value  = "Albahaca"
$ws.Cells.Item(18, 8).Value  = "Sin especificar"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 140
$ws.Cells.Item(18, 11).Value = 4500
$ws.Cells.Item(18, 12).Value = 5000
$ws.Cells.Item(18, 13).Value = 4714
$ws.Cells.Item(18, 14).Value = "`$/docena de matas"
$ws.Cells.Item(18, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 16).Value = 786
$ws.Cells.Item(18, 17).Value = 6
$ws.Cells.Item(18, 18).Value = "Hortaliza"
